$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# MISC BOM sheet: add a new "PJRC" table (rows 11-14) below the existing
# KEYCHRON (rows 1-5) / AMAZON (rows 6-9) tables.
# ---------------------------------------------------------------------------

# Copy formatting (incl. number formats / borders) from the KEYCHRON title +
# header + first data row down into the new title/header/first-data rows,
# then copy the data-row formatting once more for the second new data row.
$ws2.Range("A1:E3").Copy()
$ws2.Range("A11:E13").PasteSpecial(-4122)
$ws2.Range("A3:E3").Copy()
$ws2.Range("A14:E14").PasteSpecial(-4122)

# Merge the new title row, same as the other two tables.
$ws2.Range("A11:E11").Merge()

# Match the plain formatting used in column F next to the other two tables.
$ws2.Range("F12:F14").HorizontalAlignment = -4131

# Title & header row heights.
$ws2.Range("A11").EntireRow.RowHeight = 26.25
$ws2.Range("A12").EntireRow.RowHeight = 15.75

# Title row.
$ws2.Range("A11").Value = "PJRC"

# Header row.
$ws2.Range("A12").Value = "ITEM"
$ws2.Range("B12").Value = "QUANTITY"
$ws2.Range("C12").Value = "UNIT PRICE"
$ws2.Range("D12").Value = "EXTENDED PRICE"
$ws2.Range("E12").Value = "LINK"

# Data row 1: Teensy 4.1
$ws2.Range("A13").Value = "Teensy 4.1"
$ws2.Range("B13").Value = 1
$ws2.Range("C13").Value = 26.85
$ws2.Range("D13").Value = 26.85
$ws2.Range("E13").Value = "http://bit.ly/3ibYI2y"

# Data row 2: PSRAM Chip
$ws2.Range("A14").Value = "PSRAM Chip"
$ws2.Range("B14").Value = 2
$ws2.Range("C14").Value = 1.2
$ws2.Range("D14").Value = 2.4
$ws2.Range("E14").Value = "http://bit.ly/3nQR6nC"

# Hyperlinks for the new link cells. Hyperlinks.Add resets the cell style to
# a generic hyperlink style, so reapply the "Hyperlink" cell-style formatting
# (copied from the existing E3 hyperlink cell) straight afterwards.
$ws2.Hyperlinks.Add($ws2.Range("E13"), "http://bit.ly/3ibYI2y")
$ws2.Range("E3").Copy()
$ws2.Range("E13").PasteSpecial(-4122)

$ws2.Hyperlinks.Add($ws2.Range("E14"), "http://bit.ly/3nQR6nC")
$ws2.Range("E3").Copy()
$ws2.Range("E14").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# View / selection changes: MISC BOM becomes the active tab, selections move.
# ---------------------------------------------------------------------------
$ws1.Range("A53:H53").Select()
$ws2.Activate()
$ws2.Range("E38").Select()
